$wb = $excel.ActiveWorkbook

# ALC row 96
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1082.8
$ws.Range("I96").Value = 1253.5
$ws.Range("J96").Value = 400
$ws.Range("K96").Value = 3760.5
$ws.Range("L96").Value = 1200
$ws.Range("M96").Value = -2387.5
$ws.Range("N96").Value = -3946

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 593000.7
$ws.Range("I98").Value = 662031
$ws.Range("J98").Value = 6243
$ws.Range("K98").Value = 662031
$ws.Range("L98").Value = 6243
$ws.Range("M98").Value = -660533
$ws.Range("N98").Value = -9239

# ALC row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 551.3461
$ws.Range("I103").Value = 506.5
$ws.Range("K103").Value = 1519.5
$ws.Range("M103").Value = -933.5

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 593000.7
$ws.Range("I122").Value = 662031
$ws.Range("J122").Value = 6243
$ws.Range("K122").Value = 1986093
$ws.Range("L122").Value = 18729
$ws.Range("M122").Value = -1983643
$ws.Range("N122").Value = -23629

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 25430.121
$ws.Range("I132").Value = 27385.475
$ws.Range("J132").Value = 662.3333
$ws.Range("K132").Value = 82156.42499999999
$ws.Range("L132").Value = 1986.9999
$ws.Range("M132").Value = -79626.42499999999
$ws.Range("N132").Value = -7046.9999

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3241.94
$ws.Range("I61").Value = 2575.7144
$ws.Range("J61").Value = 6739.625
$ws.Range("K61").Value = 2575.7144
$ws.Range("L61").Value = 6739.625
$ws.Range("M61").Value = -2363.7144
$ws.Range("N61").Value = -7163.625

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2422.7673
$ws.Range("I132").Value = 2012.9459
$ws.Range("K132").Value = 6038.8377
$ws.Range("M132").Value = -3508.8377

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3241.94
$ws.Range("I136").Value = 2575.7144
$ws.Range("J136").Value = 6739.625
$ws.Range("K136").Value = 7727.1432
$ws.Range("L136").Value = 20218.875
$ws.Range("M136").Value = -5177.1432
$ws.Range("N136").Value = -25318.875

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1783.3334
$ws.Range("I86").Value = 1925
$ws.Range("K86").Value = 1925
$ws.Range("M86").Value = -802

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1783.3334
$ws.Range("I89").Value = 1925
$ws.Range("K89").Value = 9625
$ws.Range("M89").Value = -4009

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 281245.6
$ws.Range("I105").Value = 3434.95
$ws.Range("J105").Value = 628508.9
$ws.Range("K105").Value = 3434.95
$ws.Range("L105").Value = 628508.9
$ws.Range("M105").Value = -1687.95
$ws.Range("N105").Value = -632002.9

# CRP row 20
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

# CRP row 30
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2385.0938
$ws.Range("I58").Value = 1106.8
$ws.Range("J58").Value = 3513
$ws.Range("K58").Value = 1106.8
$ws.Range("L58").Value = 3513
$ws.Range("M58").Value = -903.8
$ws.Range("N58").Value = -3919

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1010.1667
$ws.Range("I122").Value = 1092.2
$ws.Range("J122").Value = 600
$ws.Range("K122").Value = 3276.6
$ws.Range("L122").Value = 1800
$ws.Range("M122").Value = -826.6000000000004
$ws.Range("N122").Value = -6700

# CRP row 125
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H125").Value = 60000
$ws.Range("J125").Value = 60000
$ws.Range("L125").Value = 60000
$ws.Range("N125").Value = -64920

# CRP row 127
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

# CRP row 128
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

# CRP row 129
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999

# CRP row 131
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H131").Value = 33500
$ws.Range("J131").Value = 33500
$ws.Range("L131").Value = 33500
$ws.Range("N131").Value = -43580

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2385.0938
$ws.Range("I136").Value = 1106.8
$ws.Range("J136").Value = 3513
$ws.Range("K136").Value = 3320.4
$ws.Range("L136").Value = 10539
$ws.Range("M136").Value = -770.3999999999996
$ws.Range("N136").Value = -15639

# CUL row 119
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 568.8570999999999
$ws.Range("I119").Value = 568.8570999999999
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 1706.5713
$ws.Range("L119").Value = 0
$ws.Range("M119").Value = 3131.4287
$ws.Range("N119").ClearContents()

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1396.2941
$ws.Range("I122").Value = 299.6
$ws.Range("J122").Value = 1853.25
$ws.Range("K122").Value = 2696.4
$ws.Range("L122").Value = 16679.25
$ws.Range("M122").Value = -246.4000000000001
$ws.Range("N122").Value = -21579.25

# CUL row 127
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 1256.3684
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 1256.3684
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 3769.1052
$ws.Range("M127").ClearContents()
$ws.Range("N127").Value = -13689.1052

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 7093639
$ws.Range("I131").Value = 443.6
$ws.Range("K131").Value = 1330.8
$ws.Range("M131").Value = 3709.2

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1769.56
$ws.Range("I122").Value = 1526.3684
$ws.Range("J122").Value = 2539.6667
$ws.Range("K122").Value = 4579.1052
$ws.Range("L122").Value = 7619.000100000001
$ws.Range("M122").Value = -2129.1052
$ws.Range("N122").Value = -12519.0001

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2646.842
$ws.Range("I126").Value = 2248.5715
$ws.Range("J126").Value = 2879.1667
$ws.Range("K126").Value = 6745.7145
$ws.Range("L126").Value = 8637.500100000001
$ws.Range("M126").Value = -4275.7145
$ws.Range("N126").Value = -13577.5001

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2503.4583
$ws.Range("I132").Value = 2146.9048
$ws.Range("J132").Value = 4999.3335
$ws.Range("K132").Value = 6440.714399999999
$ws.Range("L132").Value = 14998.0005
$ws.Range("M132").Value = -3910.714399999999
$ws.Range("N132").Value = -20058.0005

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2746.963
$ws.Range("I7").Value = 1558.5
$ws.Range("J7").Value = 3247.3684
$ws.Range("K7").Value = 1558.5
$ws.Range("L7").Value = 3247.3684
$ws.Range("M7").Value = -1446.5
$ws.Range("N7").Value = -3471.3684

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2746.963
$ws.Range("I126").Value = 1558.5
$ws.Range("J126").Value = 3247.3684
$ws.Range("K126").Value = 4675.5
$ws.Range("L126").Value = 9742.1052
$ws.Range("M126").Value = -2205.5
$ws.Range("N126").Value = -14682.1052

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4934.517
$ws.Range("I132").Value = 3164.3635
$ws.Range("K132").Value = 9493.0905
$ws.Range("M132").Value = -6963.0905

# WVR row 64
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 29057
$ws.Range("J64").Value = 29057
$ws.Range("L64").Value = 29057
$ws.Range("N64").Value = -29553

# WVR row 67
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H67").Value = 29057
$ws.Range("J67").Value = 29057
$ws.Range("L67").Value = 29057
$ws.Range("N67").Value = -30773

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 48522.285
$ws.Range("I126").Value = 91388.45
$ws.Range("J126").Value = 1369.5
$ws.Range("K126").Value = 274165.35
$ws.Range("L126").Value = 4108.5
$ws.Range("M126").Value = -271695.35
$ws.Range("N126").Value = -9048.5

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2531.7761
$ws.Range("I132").Value = 2704.62
$ws.Range("K132").Value = 8113.86
$ws.Range("M132").Value = -5583.86
